$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 9 with the "bilibili" enemy entry, mirroring the layout
# used by the existing rows (id, name, path, cost, spawnY, speed, score,
# refreshTime, refreshAccelerate).
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "bilibili"
$ws.Range("C9").Value = "res://scene/enemy/BiliBili.tscn"
$ws.Range("D9").Value = 50
$ws.Range("E9").Value = 450
$ws.Range("F9").Value = 1000
$ws.Range("G9").Value = 40
$ws.Range("H9").Value = 15
$ws.Range("I9").Value = 2

# Match the style of the other data rows (centered alignment).
$ws.Range("A9:I9").HorizontalAlignment = -4108
$ws.Range("A9:I9").VerticalAlignment = -4108

# Update the active selection as recorded in the saved workbook.
$ws.Range("H14").Select()
